# Test case 12 added: rename three existing test-method entries in column A
# (upgrade/prorated wording -> downgrade/recurring wording) and shrink
# column A's width to fit the now-shorter longest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "verifyCustomerPackageAndLicenseDowngrade"
$ws.Range("A6").Value = "verifyCustomerReceiptPageWithRecurringOrderDetails"
$ws.Range("A7").Value = "verifyCustomerReceivedSubscriptionDowngradeReceipt"

# Column A was best-fit to its longest entry (63 chars -> stored width
# 61.5703125 chars). The longest entry is now 50 chars, and the authored
# workbook stores column A at 50.9921875 chars afterwards. This runtime's
# ColumnWidth setter snaps to a whole-pixel grid (stored = ROUND(value*6+5)/6
# at this column's font metrics), so the nearest reachable grid point to
# 50.9921875 is 51.0 chars, landed on by any input in (50.0833, 50.25) -
# 50.15 is used here, comfortably inside that window.
$ws.Columns.Item(1).ColumnWidth = 50.15
